$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 306
$ws.Range("F3").Value = 1079
$ws.Range("F4").Value = 1236
$ws.Range("F5").Value = 1109
$ws.Range("F6").Value = 3319
$ws.Range("F8").Value = 49
$ws.Range("F9").Value = 1163
$ws.Range("F11").Value = 571
$ws.Range("F14").Value = 128
$ws.Range("F15").Value = 645
$ws.Range("F16").Value = 1653
$ws.Range("F17").Value = 1653
$ws.Range("F19").Value = 318
$ws.Range("F20").Value = 13
$ws.Range("F21").Value = 35
$ws.Range("F22").Value = 620
$ws.Range("F23").Value = 368
$ws.Range("F25").Value = 626
$ws.Range("F26").Value = 76927
$ws.Range("F27").Value = 76928
$ws.Range("F29").Value = 647
$ws.Range("F30").Value = 33235
$ws.Range("F31").Value = 33235
$ws.Range("F32").Value = 463
$ws.Range("F33").Value = 12
$ws.Range("F34").Value = 8
$ws.Range("F36").Value = 9
$ws.Range("F38").Value = 249
$ws.Range("F40").Value = 522
$ws.Range("F41").Value = 1161
$ws.Range("F42").Value = 5390
$ws.Range("F43").Value = 727
$ws.Range("F44").Value = 432
$ws.Range("F47").Value = 343

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F15").Value = 995
$ws.Range("F17").Value = 68
$ws.Range("F18").Value = 397
$ws.Range("F25").Value = 763
$ws.Range("F35").Value = 1477
$ws.Range("F38").Value = 99
$ws.Range("F39").Value = 99
$ws.Range("F46").Value = 22
$ws.Range("F47").Value = 39

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 543
$ws.Range("F6").Value = 560

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 306
$ws.Range("F5").Value = 543
$ws.Range("F6").Value = 1079
$ws.Range("F7").Value = 1236
$ws.Range("F9").Value = 1109
$ws.Range("F11").Value = 3319
$ws.Range("F14").Value = 49
$ws.Range("F15").Value = 1163
$ws.Range("F19").Value = 560
$ws.Range("F21").Value = 571
$ws.Range("F24").Value = 1653
$ws.Range("F25").Value = 1653
$ws.Range("F27").Value = 318
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 35
$ws.Range("F31").Value = 620
$ws.Range("F32").Value = 368
$ws.Range("F33").Value = 626
$ws.Range("F34").Value = 76930
$ws.Range("F35").Value = 647
$ws.Range("F36").Value = 33235
$ws.Range("F37").Value = 463
$ws.Range("F38").Value = 12
$ws.Range("F39").Value = 8
$ws.Range("F41").Value = 8
$ws.Range("F43").Value = 249
$ws.Range("F45").Value = 522
$ws.Range("F47").Value = 5390
$ws.Range("F49").Value = 99
